# Sample Project / Main.xlsx — "Rules" sheet
# B11 currently holds the text "R40" (last rule-row label). The save
# being replayed here overwrites that cell with the literal text "1"
# (a text value, not the number 1 — it must stay a string so it keeps
# rendering/comparing like the other rule-row labels in column B).
#
# A plain  $ws.Range("B11").Value = "1"  would let Excel's normal
# type-inference kick in and store it as the NUMBER 1 (and would also
# mint a new cell style if we tried to force text via NumberFormat).
# To get a genuine text cell without touching the cell's style, build
# the text with a formula (TEXT() always returns a string) and then
# collapse the formula down to its static cached value with a
# Copy / Paste-Special-Values round trip — exactly what a user does
# in Excel to "freeze" a formula result as a literal value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")

# TEXT(...) forces a string result, so the formula cell is t="str".
$target.Formula = '=TEXT(1,"0")'

# Copy it and paste only the value back onto itself: this turns the
# formula into a literal shared-string "1" while leaving the cell's
# existing style/number format completely untouched.
$target.Copy()
$target.PasteSpecial(-4163)   # xlPasteValues

$target.Select()
